$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds numeric-looking text (e.g. "0.7123",
# "29.308.11") that must stay plain text, matching the inlineStr cells
# in the source workbook. Temporarily force Text format on the whole
# Price column so assigning .Value does not get auto-coerced into a
# number, then restore the default "Normal" style so no stray style
# index is left behind on the cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '29.308.11'
$ws.Range('E2').Value = '  -0.05%  '
$ws.Range('D3').Value = '1.874.76'
$ws.Range('E3').Value = '  +0.03%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '0.7123'
$ws.Range('E5').Value = '  -0.05%  '
$ws.Range('D6').Value = '242.49'
$ws.Range('E6').Value = '  +0.29%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '0.08011'
$ws.Range('E8').Value = '  +2.94%  '
$ws.Range('D9').Value = '0.3154'
$ws.Range('E9').Value = '  +1.35%  '
$ws.Range('D10').Value = '24.97'
$ws.Range('E10').Value = '  -0.70%  '
$ws.Range('D11').Value = '0.08226'
$ws.Range('E11').Value = '  -2.14%  '
$ws.Range('D12').Value = '1.884.92'
$ws.Range('E12').Value = '  +0.58%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '5.247'
$ws.Range('E13').Value = '  -0.02%  '
$ws.Range('B14').Value = 'Litecoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D14').Value = '94.77'
$ws.Range('E14').Value = '  +3.92%  '
$ws.Range('D15').Value = '0.7113'
$ws.Range('E15').Value = '  -0.13%  '
$ws.Range('D16').Value = '6.417'
$ws.Range('E16').Value = '  +5.79%  '
$ws.Range('D17').Value = '0.000008549'
$ws.Range('E17').Value = '  +4.27%  '
$ws.Range('D18').Value = '29.323.99'
$ws.Range('E18').Value = '  -0.03%  '
$ws.Range('D19').Value = '243.67'
$ws.Range('E19').Value = '  +1.38%  '
$ws.Range('D20').Value = '2.146.17'
$ws.Range('E20').Value = '  +1.24%  '
$ws.Range('D21').Value = '13.24'
$ws.Range('E21').Value = '  +0.16%  '
$ws.Range('E22').Value = '  +0.12%  '
$ws.Range('D23').Value = '7.771'
$ws.Range('E23').Value = '  +0.05%  '
$ws.Range('D24').Value = '1.001'
$ws.Range('E24').Value = '  -0.07%  '
$ws.Range('D25').Value = '0.1560'
$ws.Range('E25').Value = '  -1.72%  '
$ws.Range('D26').Value = '9.040'
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('D27').Value = '162.40'
$ws.Range('E27').Value = '  -0.41%  '
$ws.Range('D28').Value = '18.52'
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('E29').Value = '  -0.55%  '
$ws.Range('D30').Value = '4.412'
$ws.Range('E30').Value = '  +0.12%  '
$ws.Range('D32').Value = '0.05377'
$ws.Range('E32').Value = '  +1.47%  '
$ws.Range('E33').Value = '  -9.52%  '
$ws.Range('D34').Value = '1.940'
$ws.Range('E34').Value = '  +0.04%  '
$ws.Range('D35').Value = '0.7628'
$ws.Range('E35').Value = '  +2.44%  '
$ws.Range('D36').Value = '1.177'
$ws.Range('E36').Value = '  +0.04%  '
$ws.Range('E37').Value = '  -0.65%  '
$ws.Range('D38').Value = '0.01875'
$ws.Range('E38').Value = '  -0.05%  '
$ws.Range('D39').Value = '1.259.14'
$ws.Range('E39').Value = '  +2.74%  '
$ws.Range('D40').Value = '2.753'
$ws.Range('E40').Value = '  +0.83%  '
$ws.Range('D41').Value = '6.476'
$ws.Range('E41').Value = '  +0.01%  '
$ws.Range('D42').Value = '0.9120'
$ws.Range('E42').Value = '  +2.96%  '
$ws.Range('D43').Value = '112.76'
$ws.Range('E43').Value = '  +3.07%  '
$ws.Range('D44').Value = '74.05'
$ws.Range('E44').Value = '  +2.12%  '
$ws.Range('D45').Value = '0.00000000133'
$ws.Range('E45').Value = '  +8.45%  '
$ws.Range('E46').Value = '  -0.04%  '
$ws.Range('D47').Value = '2.042.79'
$ws.Range('E47').Value = '  +1.36%  '
$ws.Range('D48').Value = '0.5222'
$ws.Range('E48').Value = '  +0.53%  '
$ws.Range('D49').Value = '1.798'
$ws.Range('E49').Value = '  -0.09%  '
$ws.Range('D50').Value = '9.461'
$ws.Range('E50').Value = '  +0.86%  '
$ws.Range('D51').Value = '0.4349'
$ws.Range('E51').Value = '  +0.85%  '

$ws.Range("D2:D51").Style = "Normal"
